$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

$cell1 = $table.Cell(1, 2)
$cell1.Range.Font.Bold = $true
$cell1.Range.Find.Execute("Tasa promedio de defectos", $true, $false, $false, $false, $false, $true, 1, $false, "Tasa media de defectos", 2)

$cell2 = $table.Cell(1, 3)
$cell2.Range.Font.Bold = $true
$cell2.Range.Find.Execute("Coste promedio de mantenimiento", $true, $false, $false, $false, $false, $true, 1, $false, "Promedio de coste por año de mantenimiento", 2)
